# Integrated a split function
# Append extra names (split out from elsewhere) onto a few existing
# attendee/name cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9 originally stores its text with a forced "quote prefix" (it starts
# with digits), so re-apply the leading apostrophe when rewriting it to
# keep that same text formatting.
$ws.Range("B9").Value = "'23 Mom, ryan"

$ws.Range("B11").Value = "Drew Bevington Senior, ryan, mom, trevor"

$ws.Range("B17").Value = "Trevor John, Drew, Chloe"
